# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 2004
$wsExhibit.Range("F5").Value = 317
$wsExhibit.Range("F7").Value = 91
$wsExhibit.Range("F8").Value = 2047
$wsExhibit.Range("F9").Value = 10386

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2004
$wsAll.Range("F5").Value = 317
$wsAll.Range("F8").Value = 91
$wsAll.Range("F9").Value = 2047
$wsAll.Range("F12").Value = 10386
